$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "License Information" Heading2 paragraph entirely.
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("License Information", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Expand(4)            # wdParagraph -> grab the whole paragraph incl. its mark
$rng.Delete()

# ------------------------------------------------------------------
# 2. Merge the licensing paragraph ("Perguntas de Tradução
#    (unfoldingWord) (Portuguese) is based on: ...") with the
#    following paragraph ("This PDF version is provided under the
#    same license.") by deleting the paragraph mark between them.
# ------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("is based on", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.Expand(4)            # wdParagraph -> the whole "is based on" paragraph
$mergeMark = $d.Range($rng2.End - 1, $rng2.End)
$mergeMark.Delete()

# ------------------------------------------------------------------
# 3. Replace the whole contents of the (now merged) paragraph with
#    the new resource-license text.
# ------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("is based on", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng3.Expand(4)            # wdParagraph
$body = $d.Range($rng3.Start, $rng3.End - 1)
$body.Text = ""

$insertionPoint = $d.Range($body.Start, $body.Start)
$insertionPoint.InsertAfter("unfoldingWord® Translation Questions")
$insertionPoint.Collapse(0)
$insertionPoint.InsertAfter(" © 2022 unfoldingWord. Released under CC BY-SA 4.0 license. ")
$insertionPoint.Collapse(0)
$insertionPoint.InsertAfter("unfoldingWord® Translation Questions")
$insertionPoint.Collapse(0)
$insertionPoint.InsertAfter(" has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文) from ")
$insertionPoint.Collapse(0)
$insertionPoint.InsertAfter("unfoldingWord® Translation Questions")
$insertionPoint.Collapse(0)
$insertionPoint.InsertAfter(" © 2022 unfoldingWord. Released under CC BY-SA 4.0 license by Mission Mutual")

# ------------------------------------------------------------------
# 4. Bold only the first "unfoldingWord® Translation Questions" run.
# ------------------------------------------------------------------
$boldRange = $d.Range($body.Start, $body.End)
$boldFound = $boldRange.Find.Execute("unfoldingWord® Translation Questions", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boldRange.Bold = 1

Write-Output "done"
